# Correction in SA algorithm fitness values (column C) for run_16 log.
# Column C ("Fitness") values are updated in 5 contiguous blocks of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row start, row end, new fitness value
$blocks = @(
    @(2,   9,   8012),
    @(10,  81,  8010),
    @(82,  101, 7867),
    @(102, 113, 7865),
    @(114, 252, 7293)
)

foreach ($block in $blocks) {
    $startRow = $block[0]
    $endRow   = $block[1]
    $value    = $block[2]
    $range = $ws.Range("C$startRow`:C$endRow")
    $range.Value = $value
}
